$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 84 (shifts existing rows 84..166 down to 85..167)
$ws.Rows.Item(84).Insert()

$ws.Cells.Item(84, 1).Value = 6
$ws.Cells.Item(84, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(84, 3).Value = "Metropolitana"
$ws.Cells.Item(84, 4).Value = 44601
$ws.Cells.Item(84, 5).Value = 13
$ws.Cells.Item(84, 6).Value = 100112029
$ws.Cells.Item(84, 7).Value = "Orégano"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 45
$ws.Cells.Item(84, 11).Value = 15000
$ws.Cells.Item(84, 12).Value = 15000
$ws.Cells.Item(84, 13).Value = 15000
$ws.Cells.Item(84, 14).Value = "$/docena de atados"
$ws.Cells.Item(84, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(84, 16).Value = 5000
$ws.Cells.Item(84, 17).Value = 3
$ws.Cells.Item(84, 18).Value = "Hortaliza"
